# xlwings tests/tables.xlsx fixture update:
# - rename Sheet1 -> template
# - shrink/move Table715 (table7) and Table1116 (table8) on the "expected" sheet
#   using ClearContents (not Delete) on the vacated range, matching the
#   table.update() fix described in the commit message.

$wb = $excel.ActiveWorkbook
$wsTemplate = $wb.Worksheets.Item(1)
$wsExpected = $wb.Worksheets.Item(2)

$wsTemplate.Name = "template"

$loTable715 = $wsExpected.ListObjects.Item("Table715")
$loTable1116 = $wsExpected.ListObjects.Item("Table1116")

$oldRange715 = $loTable715.Range
$oldRange1116 = $loTable1116.Range

# Snapshot values before either range is touched - Table1116's old range
# (A37:D41) overlaps Table715's new range (A39:E43), so read everything
# first, then clear, then write.
$vals715 = $oldRange715.Value()
$vals1116 = $oldRange1116.Value()

$oldRange715.ClearContents()
$oldRange1116.ClearContents()

$newRange715 = $wsExpected.Range("A39:E43")
$newRange715.Value = $vals715
$loTable715.Resize($newRange715)

$newRange1116 = $wsExpected.Range("A52:D56")
$newRange1116.Value = $vals1116
$loTable1116.Resize($newRange1116)

$wsExpected.Range("A23").Select()
$wsTemplate.Activate()
